{"js": "// Locate the target paragraph: \"The null hypothesis would be that there is no\n// correlation between happiness levels and the economic performance , this+\"\n// and replace it with the completed Introduction text, split across three\n// paragraphs (null hypothesis / alternative hypothesis / assumption).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"The null hypothesis would be that there is no correlation\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the target 'null hypothesis' paragraph.\");\n}\n\n// Replace the whole paragraph's text with the completed first sentence.\ntarget.insertText(\n  \"The null hypothesis  is that there is no correlation between the level of happiness and the economic performance of a country , suggesting that an increase in a countries or an individual\\u2019s income would not correlate with an increase in happiness.\",\n  \"Replace\"\n);\n\n// Insert the \"alternative hypothesis\" paragraph right after it.\nconst altPara = target.insertParagraph(\n  \"The alternative hypothesis is that there is a correlation between the level of happiness and the economic performance of a country, suggesting that an increase in income leads to an increase in happiness.\",\n  \"After\"\n);\n\n// Insert the \"assumption\" paragraph right after that.\naltPara.insertParagraph(\n  \"Our assumption before performing this study is that the alternative hypothesis is accepted and the null hypothesis is rejected, this may be proven or disproven by the study.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Locate the target paragraph: \"The null hypothesis would be that there is no\n# correlation between happiness levels and the economic performance , this+\"\n# and replace it with the completed Introduction text, split across three\n# paragraphs (null hypothesis / alternative hypothesis / assumption).\n$d = $word.ActiveDocument\n\n$target = $null\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"*The null hypothesis would be that there is no correlation*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate the target 'null hypothesis' paragraph.\"\n}\n\n# Replace the whole paragraph's text (use a fresh Range over the same span so\n# the trailing run gets fully overwritten instead of left behind).\n$r = $d.Range($target.Range.Start, $target.Range.End)\n$r.Text = \"The null hypothesis  is that there is no correlation between the level of happiness and the economic performance of a country , suggesting that an increase in a countries or an individual\" + [char]0x2019 + \"s income would not correlate with an increase in happiness.\"\n\n# Insert the \"alternative hypothesis\" paragraph right after it.\n$r.Collapse(0)\n$r.InsertAfter(\"`rThe alternative hypothesis is that there is a correlation between the level of happiness and the economic performance of a country, suggesting that an increase in income leads to an increase in happiness.\")\n\n# Insert the \"assumption\" paragraph right after that.\n$r.Collapse(0)\n$r.InsertAfter(\"`rOur assumption before performing this study is that the alternative hypothesis is accepted and the null hypothesis is rejected, this may be proven or disproven by the study.\")\n"}
